# analysis for pilot 2.1-2.3, power analysis for final data collection
#
# Adds subject FYP_20190410_01 (row 24) to the "pilot" sheet, scrolls the
# frozen header pane down so the new row is visible, and updates the
# remembered selection on the "conditions" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "conditions" sheet: remembered selection moves from D17 to E14.
#    Do this first so the workbook ends up with "pilot" as the final
#    active sheet/tab (matching the original file).
# ---------------------------------------------------------------------
$wsCond = $wb.Worksheets.Item("conditions")
$wsCond.Activate()
$wsCond.Range("E14").Select()

# ---------------------------------------------------------------------
# 2) "pilot" sheet: new data row for subject 23 (bing_id ba815),
#    collected 2019-04-10.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("pilot")
$ws.Activate()

# Write E24 before A24 so the new shared-string entries are interned in
# the same order as the target file ("ba815" = index 120, then
# "FYP_20190410_01" = index 121).
$ws.Range("E24").Value = "ba815"
$ws.Range("A24").Value = "FYP_20190410_01"

$ws.Range("B24").Value = 23

# C24 (test_date) / D24 (location) / O24 (age) reuse the number formats
# already used by row 23 -- copy those formats over first, then set the
# values, so no redundant numFmt entries get created.
$ws.Range("C23").Copy()
$ws.Range("C24").PasteSpecial(-4122)
$ws.Range("C24").Value = 43565

$ws.Range("D23").Copy()
$ws.Range("D24").PasteSpecial(-4122)
$ws.Range("D24").Value = "Bing"

$ws.Range("O23").Copy()
$ws.Range("O24").PasteSpecial(-4122)
$ws.Range("O24").Value = 5.4712328767123291

$ws.Range("K24").Value = "no"
$ws.Range("L24").Value = "no"
$ws.Range("M24").Value = "no"
$ws.Range("N24").Value = "male"
$ws.Range("P24").Value = "yes"
$ws.Range("Q24").Value = "One of them is bigger"
$ws.Range("R24").Value = "The yellow bucket"
$ws.Range("S24").Value = "One of them is bigger"
$ws.Range("T24").Value = "The green bucket"
$ws.Range("U24").Value = "Yellow-Ball"
$ws.Range("W24").Value = "Green-Ball"
$ws.Range("AA24").Value = "Because of the size of the buckets in her classroom"
$ws.Range("AE24").Value = "Maybe Green-Ball"
$ws.Range("AG24").Value = "Okay"
$ws.Range("AK24").Value = "For sure Yellow-Ball"

# ---------------------------------------------------------------------
# 3) Scroll the frozen header pane down so row 19 is the first visible
#    row beneath the header (keeping the existing 1-row freeze), then
#    restore the selection to A24.
# ---------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 19
$ws.Range("A24").Select()

Write-Host "done"
